$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.517.55"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "3.064.51"
$ws.Range("E3").Value = "  -3.71%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'588.06"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "'154.12"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.535"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "3.059.39"
$ws.Range("E9").Value = "  -3.99%  "
$ws.Range("E10").Value = "  -4.43%  "
$ws.Range("D11").Value = "'5.83"
$ws.Range("E11").Value = "  -2.16%  "
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").Value = "'36.70"
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("E14").Value = "  -4.94%  "
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").Value = "3.571.70"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("D17").Value = "63.522.53"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "'7.10"
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("D19").Value = "3.066.29"
$ws.Range("E19").Value = "  -3.59%  "
$ws.Range("D20").Value = "'468.34"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "'14.22"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").Value = "'0.701"
$ws.Range("E22").Value = "  -5.61%  "
$ws.Range("D23").Value = "'7.44"
$ws.Range("E23").Value = "  -3.64%  "
$ws.Range("D24").Value = "'2.41"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "'80.30"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").Value = "'12.71"
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("D27").Value = "'10.33"
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'7.35"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("D32").Value = "'2.13"
$ws.Range("E32").Value = "  -6.22%  "
$ws.Range("D33").Value = "'27.01"
$ws.Range("E33").Value = "  -5.07%  "
$ws.Range("E34").Value = "  -6.09%  "
$ws.Range("D35").Value = "0.0₃0815"
$ws.Range("E35").Value = "  -5.64%  "
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").Value = "'5.95"
$ws.Range("E37").Value = "  -4.71%  "
$ws.Range("D38").Value = "'3.24"
$ws.Range("E38").Value = "  -2.84%  "
$ws.Range("D39").Value = "'2.19"
$ws.Range("E39").Value = "  -5.71%  "
$ws.Range("D40").Value = "'50.41"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").Value = "'9.15"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("D42").Value = "'435.39"
$ws.Range("E42").Value = "  -8.38%  "
$ws.Range("D43").Value = "'0.285"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("D44").Value = "'40.51"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").Value = "'0.110"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "'0.0357"
$ws.Range("E46").Value = "  -5.49%  "
$ws.Range("D47").Value = "2.790.01"
$ws.Range("E47").Value = "  -4.57%  "
$ws.Range("D48").Value = "'128.86"
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'24.85"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").Value = "'2.20"
$ws.Range("E51").Value = "  -3.48%  "
